# Added institutional distinctiveness in IQAC
# Adds a new staff row (THAMARAI SELVAN P, Lab Instructor) to the staff
# data sheet, and normalizes the formatting of the existing "extra"
# rows (19-23) that previously carried a redundant "no-op" style on
# columns B/C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Strip the redundant (visually identical to default) style that
#        was applied to columns B and C of rows 19-23. Doing this via
#        Borders removal collapses the cell back onto the default style
#        index instead of leaving a duplicate style behind.
$ws.Range("B19:C23").Borders.LineStyle = -4142

# --- 2. Row 19 no longer needs an explicit custom row height -- let it
#        fall back to the sheet default (14.4) like rows 20-23 already do.
$ws.Rows.Item(19).AutoFit() | Out-Null

# --- 3. Append the new staff member as row 24. Set the cell values
#        first, then paste the matching formatting from the equivalent
#        cells one row up so the bordered "name" style carries over
#        exactly instead of being dropped by the later value write.
$ws.Range("A24").Value = "THAMARAI SELVAN P"
$ws.Range("B24").Value = "Lab Instructor"
$ws.Range("C24").Value = "/static/images/profile_photos/001/VEC-001-05-1.webp"
$ws.Range("J24").Value = "VEC-001-05-1"

$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A24").Copy()
$ws.Range("B24").PasteSpecial(-4122)

$ws.Range("D19").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("F19").Copy()
$ws.Range("F24").PasteSpecial(-4122)
$ws.Range("G19").Copy()
$ws.Range("G24").PasteSpecial(-4122)
$ws.Range("H19").Copy()
$ws.Range("H24").PasteSpecial(-4122)
$ws.Range("I19").Copy()
$ws.Range("I24").PasteSpecial(-4122)
$ws.Range("J19").Copy()
$ws.Range("J24").PasteSpecial(-4122)

$ws.Rows.Item(24).RowHeight = 18.75

Write-Host "Added THAMARAI SELVAN P as row 24"
